{"js": "// Remove the \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block (and the blank\n// paragraph immediately preceding it) that used to trail the bibliography\n// text in the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"Artigos de revistas especializadas...\").\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Artigos de revistas especializadas\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the anchor paragraph.\");\n}\n\n// The three paragraphs that immediately follow the anchor are the blank\n// separator paragraph, the \"Ver no Jupiter...\" line and the \"\u00a9 2020 ...\"\n// line \u2014 all of which should be removed.\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i < items.length && toDelete.length < 3; i++) {\n  const text = items[i].text;\n  if (\n    text === \"\" ||\n    text.indexOf(\"Ver no Jupiter\") !== -1 ||\n    text.indexOf(\"\\u00A9 2020\") !== -1 ||\n    text.indexOf(\"Contact: luizeleno@usp.br\") !== -1\n  ) {\n    toDelete.push(items[i]);\n  } else {\n    break;\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block (and the blank\n# paragraph immediately preceding it) that used to trail the bibliography\n# text in the document body.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($text -like \"*Artigos de revistas especializadas*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the anchor paragraph.\"\n}\n\n# The three paragraphs that immediately follow the anchor are the blank\n# separator paragraph, the \"Ver no Jupiter...\" line and the \"\u00a9 2020 ...\"\n# line \u2014 all of which should be removed. Repeatedly deleting the paragraph\n# right after the anchor keeps the index valid as the collection shrinks.\nfor ($k = 0; $k -lt 3; $k++) {\n    $target = $d.Paragraphs.Item($anchorIndex + 1).Range\n    $target.Delete()\n}\n"}
